# Updates the cryptos list to refresh prices/volume percentages and
# reorder a couple of rows (WrappedEther/Polygon and Kaspa/RocketPoolETH swapped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "34.572.74"
$ws.Range("E2").Value = "  -0.23%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "1.814.26"
$ws.Range("E3").Value = "  -0.13%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.04%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "228.81"
$ws.Range("E5").Value = "  +0.22%  "

# --- Row 6 (XRP) ---
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  +8.32%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  +0.00%  "

# --- Row 8 (Solana) ---
$ws.Range("D8").Value = "36.69"
$ws.Range("E8").Value = "  +5.11%  "

# --- Row 9 (Cardano) ---
$ws.Range("E9").Value = "  +0.37%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("D10").Value = "0.0703"
$ws.Range("E10").Value = "  +1.04%  "

# --- Row 11 (TRON) ---
$ws.Range("D11").Value = "0.0970"
$ws.Range("E11").Value = "  +1.95%  "

# --- Row 12 (WrappedliquidstakedEther2.0) ---
$ws.Range("D12").Value = "2.076.34"
$ws.Range("E12").Value = "  -0.14%  "

# --- Row 13 (Chainlink) ---
$ws.Range("D13").Value = "11.59"
$ws.Range("E13").Value = "  +2.12%  "

# --- Row 14 (was WrappedEther, now Polygon) ---
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.657"
$ws.Range("E14").Value = "  +2.11%  "

# --- Row 15 (was Polygon, now WrappedEther) ---
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.821.26"
$ws.Range("E15").Value = "  -0.09%  "

# --- Row 16 (Polkadot) ---
$ws.Range("D16").Value = "4.50"
$ws.Range("E16").Value = "  +4.07%  "

# --- Row 17 (WrappedBTC) ---
$ws.Range("D17").Value = "34.548.72"
$ws.Range("E17").Value = "  -0.43%  "

# --- Row 18 (Litecoin) ---
$ws.Range("D18").Value = "70.45"
$ws.Range("E18").Value = "  +1.69%  "

# --- Row 19 (BitcoinCash) ---
$ws.Range("D19").Value = "247.10"
$ws.Range("E19").Value = "  -0.28%  "

# --- Row 20 (ShibaInu) ---
$ws.Range("D20").Value = "0.0`u{2083}0799"
$ws.Range("E20").Value = "  -0.49%  "

# --- Row 21 (Avalanche) ---
$ws.Range("D21").Value = "11.69"
$ws.Range("E21").Value = "  +1.43%  "

# --- Row 22 (Dai) ---
$ws.Range("E22").Value = "  +0.06%  "

# --- Row 23 (Uniswap) ---
$ws.Range("D23").Value = "4.24"
$ws.Range("E23").Value = "  +1.05%  "

# --- Row 24 (Toncoin) ---
$ws.Range("D24").Value = "2.25"
$ws.Range("E24").Value = "  +7.43%  "

# --- Row 25 (Monero) ---
$ws.Range("D25").Value = "172.48"
$ws.Range("E25").Value = "  +0.25%  "

# --- Row 26 (Cosmos) ---
$ws.Range("D26").Value = "8.04"
$ws.Range("E26").Value = "  +7.91%  "

# --- Row 27 (Stellar) ---
$ws.Range("D27").Value = "0.124"
$ws.Range("E27").Value = "  +5.40%  "

# --- Row 28 (EthereumClassic) ---
$ws.Range("D28").Value = "17.36"
$ws.Range("E28").Value = "  +3.42%  "

# --- Row 29 (BinanceUSD) ---
$ws.Range("E29").Value = "  -0.06%  "

# --- Row 30 (InternetComputer(DFINITY)) ---
$ws.Range("E30").Value = "  +1.50%  "

# --- Row 31 (Filecoin) ---
$ws.Range("E31").Value = "  +0.84%  "

# --- Row 32 (Hedera) ---
$ws.Range("D32").Value = "0.0534"
$ws.Range("E32").Value = "  +0.12%  "

# --- Row 33 (PancakeSwap) ---
$ws.Range("E33").Value = "  -0.40%  "

# --- Row 34 (LidoDAOToken) ---
$ws.Range("E34").Value = "  -1.66%  "

# --- Row 35 (Maker) ---
$ws.Range("D35").Value = "1.401.44"
$ws.Range("E35").Value = "  -1.49%  "

# --- Row 36 (ImmutableX) ---
$ws.Range("E36").Value = "  -0.57%  "

# --- Row 37 (RenderToken) ---
$ws.Range("D37").Value = "2.43"
$ws.Range("E37").Value = "  -7.94%  "

# --- Row 38 (TrustWalletToken) ---
$ws.Range("E38").Value = "  -0.07%  "

# --- Row 39 (VeChain) ---
$ws.Range("E39").Value = "  -0.08%  "

# --- Row 40 (ARBITRUM) ---
$ws.Range("E40").Value = "  +0.81%  "

# --- Row 41 (Aave) ---
$ws.Range("D41").Value = "83.24"
$ws.Range("E41").Value = "  -3.75%  "

# --- Row 42 (MXToken) ---
$ws.Range("E42").Value = "  -0.28%  "

# --- Row 43 (HuobiToken) ---
$ws.Range("E43").Value = "  +0.53%  "

# --- Row 44 (WEMIXToken) ---
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +7.84%  "

# --- Row 45 (InjectiveProtocol) ---
$ws.Range("D45").Value = "13.74"
$ws.Range("E45").Value = "  -0.18%  "

# --- Row 46 (FraxShare) ---
$ws.Range("D46").Value = "6.09"
$ws.Range("E46").Value = "  -0.95%  "

# --- Row 47 (was Kaspa, now RocketPoolETH) ---
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.976.55"
$ws.Range("E47").Value = "  -0.20%  "

# --- Row 48 (was RocketPoolETH, now Kaspa) ---
$ws.Range("B48").Value = "Kaspa"
$ws.Range("C48").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D48").Value = "0.0495"
$ws.Range("E48").Value = "  -5.84%  "

# --- Row 49 (Quant) ---
$ws.Range("D49").Value = "104.75"
$ws.Range("E49").Value = "  -1.22%  "

# --- Row 50 (PaxDollar) ---
$ws.Range("E50").Value = "  +0.03%  "

# --- Row 51 (BabyDogeCoin) ---
$ws.Range("D51").Value = "0.0`u{2086}0129"
$ws.Range("E51").Value = "  -2.02%  "
